{"js": "// This document contains a single table with 20 rows x 5 columns; every\n// cell holds one arithmetic equation (e.g. \"79-21=58\"). The edit replaces\n// each cell's text with a new equation. Several \"before\" equations are\n// duplicated across different cells but map to different \"after\" values\n// (e.g. \"90-72=18\" -> \"73+17=90\" in one cell and -> \"63-31=32\" in another),\n// so the replacement must be positional (row/column), not a global\n// find-and-replace. Assigning the full grid to Table.values does this in\n// one shot while leaving all run/paragraph formatting untouched.\nconst newValues = [\n  [\"29+16=45\", \"14+19=33\", \"11+73=84\", \"49+9=58\", \"51-5=46\"],\n  [\"82-64=18\", \"21-12=9\", \"97-7=90\", \"35-17=18\", \"6+3=9\"],\n  [\"99-10=89\", \"44+51=95\", \"92-90=2\", \"63-58=5\", \"73+17=90\"],\n  [\"85-0=85\", \"86-39=47\", \"90-82=8\", \"38+45=83\", \"5+48=53\"],\n  [\"89-74=15\", \"63+23=86\", \"58-13=45\", \"86-10=76\", \"51-2=49\"],\n  [\"63+33=96\", \"37-15=22\", \"33+28=61\", \"82-69=13\", \"83-81=2\"],\n  [\"17+16=33\", \"23-12=11\", \"69+8=77\", \"89-11=78\", \"9+48=57\"],\n  [\"29+31=60\", \"4+32=36\", \"99-18=81\", \"63-31=32\", \"37+42=79\"],\n  [\"53-0=53\", \"95-48=47\", \"39-11=28\", \"93-56=37\", \"13+1=14\"],\n  [\"62+15=77\", \"35+20=55\", \"0+81=81\", \"21+70=91\", \"31-20=11\"],\n  [\"90+4=94\", \"73-53=20\", \"27-10=17\", \"34-14=20\", \"46+11=57\"],\n  [\"58+33=91\", \"99-53=46\", \"72+25=97\", \"81-17=64\", \"8+39=47\"],\n  [\"57-29=28\", \"6+84=90\", \"73-64=9\", \"53+38=91\", \"47-29=18\"],\n  [\"83-79=4\", \"8+84=92\", \"3+96=99\", \"8+48=56\", \"54-44=10\"],\n  [\"47+29=76\", \"92-17=75\", \"94-90=4\", \"73-48=25\", \"80-49=31\"],\n  [\"53+3=56\", \"77-25=52\", \"41+18=59\", \"79-9=70\", \"71-0=71\"],\n  [\"90-12=78\", \"93-33=60\", \"8+67=75\", \"96-76=20\", \"38-18=20\"],\n  [\"90-46=44\", \"59+33=92\", \"26+56=82\", \"87-66=21\", \"11+44=55\"],\n  [\"9+15=24\", \"40+14=54\", \"20+71=91\", \"93-42=51\", \"2+68=70\"],\n  [\"86-39=47\", \"71-10=61\", \"15+28=43\", \"60-16=44\", \"67-30=37\"]\n];\n\nconst body = context.document.body;\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nif (tables.items.length === 0) {\n  throw new Error(\"Expected a table in the document body, found none.\");\n}\n\nconst table = tables.items[0];\ntable.load(\"rowCount, values\");\nawait context.sync();\n\nif (table.rowCount !== newValues.length) {\n  throw new Error(`Expected ${newValues.length} rows, found ${table.rowCount}.`);\n}\n\ntable.values = newValues;\nawait context.sync();\n", "ps1": "# Update the table of addition & subtraction equations in place.\n# The document contains exactly one table with 20 rows x 5 columns, where\n# each cell holds a single equation string (e.g. \"79-21=58\"). We replace\n# the whole grid of values positionally (row-major) via Cell.Range.Text,\n# which safely handles the cases where the same \"before\" text occurs in\n# more than one cell but must become different \"after\" text, and leaves\n# all run/paragraph formatting (fonts, size, alignment) untouched.\n$newValues = @(\n    @(\"29+16=45\", \"14+19=33\", \"11+73=84\", \"49+9=58\", \"51-5=46\"),\n    @(\"82-64=18\", \"21-12=9\", \"97-7=90\", \"35-17=18\", \"6+3=9\"),\n    @(\"99-10=89\", \"44+51=95\", \"92-90=2\", \"63-58=5\", \"73+17=90\"),\n    @(\"85-0=85\", \"86-39=47\", \"90-82=8\", \"38+45=83\", \"5+48=53\"),\n    @(\"89-74=15\", \"63+23=86\", \"58-13=45\", \"86-10=76\", \"51-2=49\"),\n    @(\"63+33=96\", \"37-15=22\", \"33+28=61\", \"82-69=13\", \"83-81=2\"),\n    @(\"17+16=33\", \"23-12=11\", \"69+8=77\", \"89-11=78\", \"9+48=57\"),\n    @(\"29+31=60\", \"4+32=36\", \"99-18=81\", \"63-31=32\", \"37+42=79\"),\n    @(\"53-0=53\", \"95-48=47\", \"39-11=28\", \"93-56=37\", \"13+1=14\"),\n    @(\"62+15=77\", \"35+20=55\", \"0+81=81\", \"21+70=91\", \"31-20=11\"),\n    @(\"90+4=94\", \"73-53=20\", \"27-10=17\", \"34-14=20\", \"46+11=57\"),\n    @(\"58+33=91\", \"99-53=46\", \"72+25=97\", \"81-17=64\", \"8+39=47\"),\n    @(\"57-29=28\", \"6+84=90\", \"73-64=9\", \"53+38=91\", \"47-29=18\"),\n    @(\"83-79=4\", \"8+84=92\", \"3+96=99\", \"8+48=56\", \"54-44=10\"),\n    @(\"47+29=76\", \"92-17=75\", \"94-90=4\", \"73-48=25\", \"80-49=31\"),\n    @(\"53+3=56\", \"77-25=52\", \"41+18=59\", \"79-9=70\", \"71-0=71\"),\n    @(\"90-12=78\", \"93-33=60\", \"8+67=75\", \"96-76=20\", \"38-18=20\"),\n    @(\"90-46=44\", \"59+33=92\", \"26+56=82\", \"87-66=21\", \"11+44=55\"),\n    @(\"9+15=24\", \"40+14=54\", \"20+71=91\", \"93-42=51\", \"2+68=70\"),\n    @(\"86-39=47\", \"71-10=61\", \"15+28=43\", \"60-16=44\", \"67-30=37\"),\n)\n\n$d = $word.ActiveDocument\n$tbl = $d.Tables.Item(1)\n\nif ($tbl.Rows.Count -ne $newValues.Count) {\n    throw \"Expected $($newValues.Count) rows, found $($tbl.Rows.Count).\"\n}\n\nfor ($r = 1; $r -le $tbl.Rows.Count; $r++) {\n    $rowValues = $newValues[$r - 1]\n    for ($c = 1; $c -le $rowValues.Count; $c++) {\n        $tbl.Cell($r, $c).Range.Text = $rowValues[$c - 1]\n    }\n}\n"}
